# Natmi following Dr Hou advice:
# A new target/sending cluster "M2" is introduced into the Rarres2/Gpr1
# ligand-receptor signalling table, expanding the Sending x Target grid
# from 3x3 to 3x4 and recomputing all the derived NATMI statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{
    A = 1;  B = 2;  C = 3;  D = 4;  E = 5;  F = 6;  G = 7;  H = 8
    I = 9;  J = 10; K = 11; L = 12; M = 13; N = 14; O = 15; P = 16
    Q = 17; R = 18; S = 19; T = 20
}

$rows = @(
    @{ Row=2; A="ECs"; B="Rarres2"; C="Gpr1"; D="ECs"; E=3; F=1; G=3.916729333333334; H=11.750188; I=0.05842616646845182; J=0.05842616646845181; K=3; L=1; M=2.75858; N=8.275739999999999; O=0.5077178793781231; P=0.507717879378123; Q=10.80461120434667; R=97.24150083912; S=0.02966400933955556; T=0.02966400933955555 }
    @{ Row=3; A="ECs"; B="Rarres2"; C="Gpr1"; D="FAPs"; E=3; F=1; G=3.916729333333334; H=11.750188; I=0.05842616646845182; J=0.05842616646845181; K=3; L=1; M=2.008641666666667; N=6.025925; O=0.3696913946416412; P=0.3696913946416412; Q=7.86730573598889; R=70.80575162390001; S=0.02159965096528665; T=0.02159965096528664 }
    @{ Row=4; A="ECs"; B="Rarres2"; C="Gpr1"; D="M2"; E=3; F=1; G=3.916729333333334; H=11.750188; I=0.05842616646845182; J=0.05842616646845181; K=1; L=0.3333333333333333; M=0.009769666666666668; N=0.029309; O=0.001798111507453522; P=0.001798111507453522; Q=0.03826514001022224; R=0.3443862600920001; S=0.0001050567622633183; T=0.0001050567622633183 }
    @{ Row=5; A="ECs"; B="Rarres2"; C="Gpr1"; D="sCs"; E=3; F=1; G=3.916729333333334; H=11.750188; I=0.05842616646845182; J=0.05842616646845181; K=3; L=1; M=0.6563016666666667; N=1.968905; O=0.1207926144727823; P=0.1207926144727823; Q=2.570555989348889; R=23.13500390414; S=0.007057449401346302; T=0.007057449401346299 }
    @{ Row=6; A="FAPs"; B="Rarres2"; C="Gpr1"; D="ECs"; E=3; F=1; G=48.53546666666667; H=145.6064; I=0.7240074597335789; J=0.7240074597335789; K=3; L=1; M=2.75858; N=8.275739999999999; O=0.5077178793781231; P=0.507717879378123; Q=133.8889676373333; R=1205.000708736; S=0.3675915321098746; T=0.3675915321098744 }
    @{ Row=7; A="FAPs"; B="Rarres2"; C="Gpr1"; D="FAPs"; E=3; F=1; G=48.53546666666667; H=145.6064; I=0.7240074597335789; J=0.7240074597335789; K=3; L=1; M=2.008641666666667; N=6.025925; O=0.3696913946416412; P=0.3696913946416412; Q=97.49036065777779; R=877.41324592; S=0.2676593275198587; T=0.2676593275198587 }
    @{ Row=8; A="FAPs"; B="Rarres2"; C="Gpr1"; D="M2"; E=3; F=1; G=48.53546666666667; H=145.6064; I=0.7240074597335789; J=0.7240074597335789; K=1; L=0.3333333333333333; M=0.009769666666666668; N=0.029309; O=0.001798111507453522; P=0.001798111507453522; Q=0.4741753308444446; R=4.2675779776; S=0.001301846144829141; T=0.001301846144829141 }
    @{ Row=9; A="FAPs"; B="Rarres2"; C="Gpr1"; D="sCs"; E=3; F=1; G=48.53546666666667; H=145.6064; I=0.7240074597335789; J=0.7240074597335789; K=3; L=1; M=0.6563016666666667; N=1.968905; O=0.1207926144727823; P=0.1207926144727823; Q=31.85390766577778; R=286.685168992; S=0.08745475395901665; T=0.08745475395901664 }
    @{ Row=10; A="sCs"; B="Rarres2"; C="Gpr1"; D="ECs"; E=3; F=1; G=14.58505066666667; H=43.755152; I=0.2175663737979692; J=0.2175663737979692; K=3; L=1; M=2.75858; N=8.275739999999999; O=0.5077178793781231; P=0.507717879378123; Q=40.23402906805333; R=362.10626161248; S=0.110462337928693; T=0.1104623379286929 }
    @{ Row=11; A="sCs"; B="Rarres2"; C="Gpr1"; D="FAPs"; E=3; F=1; G=14.58505066666667; H=43.755152; I=0.2175663737979692; J=0.2175663737979692; K=3; L=1; M=2.008641666666667; N=6.025925; O=0.3696913946416412; P=0.3696913946416412; Q=29.29614047951111; R=263.6652643156; S=0.08043241615649588; T=0.08043241615649586 }
    @{ Row=12; A="sCs"; B="Rarres2"; C="Gpr1"; D="M2"; E=3; F=1; G=14.58505066666667; H=43.755152; I=0.2175663737979692; J=0.2175663737979692; K=1; L=0.3333333333333333; M=0.009769666666666668; N=0.029309; O=0.001798111507453522; P=0.001798111507453522; Q=0.1424910833297778; R=1.282419749968; S=0.0003912086003610629; T=0.0003912086003610628 }
    @{ Row=13; A="sCs"; B="Rarres2"; C="Gpr1"; D="sCs"; E=3; F=1; G=14.58505066666667; H=43.755152; I=0.2175663737979692; J=0.2175663737979692; K=3; L=1; M=0.6563016666666667; N=1.968905; O=0.1207926144727823; P=0.1207926144727823; Q=9.572193060951111; R=86.14973754856; S=0.02628041111241934; T=0.02628041111241934 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    foreach ($col in $colIndex.Keys) {
        $ws.Cells.Item($rowNum, $colIndex[$col]).Value = $r[$col]
    }
}
